# Add two new rows of paper data (rows 4 and 5) to Sheet1, matching the
# columns already used by the existing rows 2-3 (A=Title, B=Authors,
# C=Year, D=Abstract, E=Species, M=scanner, N=neuroimaging, O=study_type).
# Columns F-L for rows 4/5 (Groups, Sex, Sample size, Age mean/sd/min/max)
# were already populated before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: White matter abnormalities / schizophrenia DTI+fMRI study
$ws.Range("A4").Value = "White matter abnormalities and brain activation in schizophrenia: A combined DTI and fMRI study"
$ws.Range("B4").Value = "Ralf G.M. Schlösser a,⁎, Igor Nenadic a,d, Gerd Wagner a, Daniel Güllmar b, Katrin von Consbruch a, Sabine Köhler a, C. Christoph Schultz a, Kathrin Koch a, Clemens Fitzek c, Paul M. Matthewsd, Jürgen R. Reichenbachb, Heinrich Saue"
$ws.Range("C4").Value = 2007
$ws.Range("D4").Value = "Diffusion tensor imaging (DTI) studies of schizophrenia have revealed white matter abnormalities in several areas of the brain. The functional impact on either psychopathology or cognition remains, however, poorly understood. Here we analysed both functional MRI (during a working memory task) and DTI data sets in 18 patients with schizophrenia and 18 controls. Firstly, DTI analyses revealed reductions of fractional anisotropy (FA) in the right medial temporal lobe adjacent to the right parahippocampal gyrus, likely to contain fibres of the inferior cingulum bundle, and in the right frontal lobe. Secondly, functional MRI revealed prefrontal, superior parietal and occipital relative hypoactivation in patients with the main effect of task. This was accounted for by reduced prefrontal activation during the encoding phase of the task, but not during maintenance or retrieval phases. Thirdly, we found a direct correlation in patients between the frontal FA reduction (but not medial temporal reductions) and fMRI activation in regions in the prefrontal and occipital cortex. Our study combining fMRI and DTI thus demonstrates altered structure-function relationships in schizophrenia. It highlights a potential relationship between anatomical changes in a frontal–temporal anatomical circuit and functional alterations in the prefrontal cortex. © 2006 Elsevier B.V. All rights reserved."
$ws.Range("E4").Value = "human"

# Row 5: MRI characteristics of substantia nigra / Parkinson's disease
$ws.Range("A5").Value = "MRI characteristics of the substantia nigra in Parkinson's disease: A combined quantitative T1 and DTI study"
$ws.Range("B5").Value = "Ricarda A. Menkea, Jan Scholza, Karla L. Millera, Sean Deonib, Saad Jbabdia, Paul M. Matthewsa,c, Mojtaba Zareia,c,⁎"
$ws.Range("C5").Value = 2009
$ws.Range("D5").Value = "The substantia nigra contains dopaminergic cells that project to the striatum and are affected by the neurodegenerative process that appears in Parkinson's disease (PD). For accurate differential diagnosis and for disease monitoring the availability of a sensitive and non-invasive biomarker for Parkinson's disease would beessential. Although there has been notable progress in studying correlates of nigral degeneration by means of magnetic resonance imaging (MRI) in the past decade, MRI and analysis techniques that allow accurate high-resolution mapping of the SN within a clinically acceptable acquisition time are still elusive. The mainpurpose of the preliminary study was to evaluate the potential role of the driven equilibrium single pulse observation of T1 (DESPOT1) method for delineation of the SN and differentiation of PD patients from healthy control subjects (n=10 in each group). We also investigated whether additional measures that can be obtained with diffusion tensor imaging (DTI) can further improve the MRI-guided discrimination between PD patients and controls. Our results show that the DESPOT1 method allows for a clear visualisation of the SN as a whole. Volumetric comparisons between ten PD patients and ten healthy subjects revealed significantly smaller volumes in patients for both the left and the right sides when the whole SN was considered. Combining SN volumetry and its connectivity with the thalamus improved the classification sensitivity to 100% and specificity to 80% for PD (discriminant function analysis with leave-one-out cross validation). Combining DESPOT1 imaging and DTI could therefore serve as a diagnostic marker for idiopathic Parkinson's disease in the future."
$ws.Range("E5").Value = "human"

# Sample size / scanner field mean (M column) for both new rows
$ws.Range("M4").Value = 3
$ws.Range("M5").Value = 3

# Neuroimaging modality tags (N column) - filled in last
$ws.Range("N4").Value = "T2,dti,resting state fMRI"
$ws.Range("N5").Value = "dti,T1-weighted"

# Study type (O column)
$ws.Range("O4").Value = "experiment"
$ws.Range("O5").Value = "experiment"

$wb.Save()
